# Auto-generated edit script applying the cryptos.xlsx data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.215.74"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "2.491.46"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.524"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0812"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "2.881.65"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "2.492.98"
$ws.Range("E16").Value = "  +2.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.847"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "47.141.75"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +15.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "246.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.87%  "
$ws.Range("E30").Value = "  +9.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0785"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.24%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.62%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.08%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0296"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("D45").Value = "1.988.10"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.11%  "
$ws.Range("E47").Value = "  -2.76%  "
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "57.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.58%  "
